$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before H (shifts creditos/periodo_idx/acred1/lugar from H:K to I:L)
$ws.Columns("H").Insert()

# New header for the inserted "reg2" column
$ws.Range("H1").Value = "reg2"

# Update the "reg1" values to include the "Resolución " prefix / full text.
# Order matters for shared-string ordering: G3 is written before G2.
$ws.Range("G3").Value = "Resolución 011706 07 de noviembre del 2019"
$ws.Range("G2").Value = "Resolución 014374 del 11 de diciembre de 2019"

# Format the new (empty) reg2 cell in row 2 like its neighbour (left/center/wrap, Arial)
# but without the border that the reg1 column has.
$ws.Range("H2").Font.Name = "Arial"
$ws.Range("H2").HorizontalAlignment = -4131
$ws.Range("H2").VerticalAlignment = -4108
$ws.Range("H2").WrapText = $true
$ws.Range("H2").Borders.LineStyle = -4142

# Widen the new reg1/reg2 columns
$ws.Range("G1:H1").EntireColumn.ColumnWidth = 54.16666666666666

# Row 2 no longer needs the extra height now that the long text moved to a wide column
$ws.Range("A2").EntireRow.RowHeight = 15

# Update the selected cell
$ws.Range("G3").Select() | Out-Null
